$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.224.50"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.788.19"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D5").Value = "'226.10"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'32.32"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'0.0691"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "'0.0947"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "2.046.54"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'11.12"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "1.787.82"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D16").Value = "34.209.64"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "'67.96"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").Value = "'245.98"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").Value = "'11.00"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "'161.94"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'16.33"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0521"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").Value = "1.441.10"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("E36").Value = "  +8.55%  "
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'0.0190"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "'82.21"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").Value = "'13.91"
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("D43").Value = "'0.921"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").Value = "'6.09"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "1.946.58"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "'105.39"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  -7.39%  "
